$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "R000"
$ws.Range("B4").Value = "Tim"
$ws.Range("C4").Value = "Test for tim"
$ws.Range("D4").Value = "2025-09-30 20:25:43"
